# Update column F (dSF) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -6
    4  = -3
    8  = -5
    9  = 6
    10 = -1
    11 = -6
    13 = -6
    14 = -1
    16 = -3
    17 = -8
    18 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
